$wb = $excel.ActiveWorkbook

# Changes common to both "展览" (Exhibition) and "全部类型" (All types) sheets.
$commonChanges = @{
    3  = 5019
    5  = 7287
    11 = 18
    12 = 4259
    13 = 1716
    14 = 96
    16 = 2862
    19 = 201
    20 = 456
    21 = 416
    22 = 443
    23 = 278
    24 = 83
    26 = 1146
    28 = 1347
    29 = 101
    30 = 570
    33 = 20
    34 = 51
    36 = 2686
    37 = 688
    38 = 36
}

# Apply to "展览" sheet (first worksheet).
$wsExhibition = $wb.Worksheets.Item(1)
foreach ($row in $commonChanges.Keys) {
    $wsExhibition.Range("F$row").Value = $commonChanges[$row]
}

# Apply to "全部类型" sheet (fourth worksheet), which also has an extra change on row 32.
$wsAllTypes = $wb.Worksheets.Item(4)
foreach ($row in $commonChanges.Keys) {
    $wsAllTypes.Range("F$row").Value = $commonChanges[$row]
}
$wsAllTypes.Range("F32").Value = 510
